$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A1").Value = 1.240880489349365
$ws.Range("B1").Value = 1.893560886383057
$ws.Range("C1").Value = 2.521998882293701
$ws.Range("D1").Value = 3.8823561668396
$ws.Range("E1").Value = 1.148517608642578
